$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column (10-nov) before column DM ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Insert a new column at DM (shifts old DM:EQ -> DN:ER)
$ws1.Range("DM1").EntireColumn.Insert()

# Header for the freshly inserted column
$ws1.Range("DM1").Value = "10-nov"

# Data rows 2-25: no data yet for 10-nov, mark with "-" like other empty days
for ($r = 2; $r -le 25; $r++) {
    $ws1.Cells.Item($r, 117).Value = "-"
}

# --- Sheet "Gaz": append two more days of data ---
$ws2 = $wb.Worksheets.Item("Gaz")

# Force the date column to stay plain text (matches existing rows) instead of
# Excel's automatic date-serial conversion.
$ws2.Range("A146").NumberFormat = "@"
$ws2.Range("A146").Value = "2025-11-08"
$ws2.Range("A146").Style = "Normal"
$ws2.Range("B146").Value = 29.755

$ws2.Range("A147").NumberFormat = "@"
$ws2.Range("A147").Value = "2025-11-09"
$ws2.Range("A147").Style = "Normal"
$ws2.Range("B147").Value = 29.755

# --- Sheet "CO2": append two more days of data ---
$ws3 = $wb.Worksheets.Item("CO2")

$ws3.Range("A146").NumberFormat = "@"
$ws3.Range("A146").Value = "2025-11-08"
$ws3.Range("A146").Style = "Normal"
$ws3.Range("B146").Value = 79.36

$ws3.Range("A147").NumberFormat = "@"
$ws3.Range("A147").Value = "2025-11-09"
$ws3.Range("A147").Style = "Normal"
$ws3.Range("B147").Value = 79.36
